$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: locate the (unique) literal text $oldText anywhere in the body and
# return a Range covering exactly that text (Find match range).
# ---------------------------------------------------------------------------
function Find-Text($doc, $oldText) {
    $rng = $doc.Content
    $ok = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        Write-Host "NOT FOUND:" $oldText
        return $null
    }
    return $doc.Range($rng.Start, $rng.End)
}

# ---------------------------------------------------------------------------
# 1) "1. What did you like the most on the website?"
#    -> "4. If applicable, describe any major challenge you may have had
#        during the tasks?"  (5 runs, all highlighted green)
# ---------------------------------------------------------------------------
$found = Find-Text $d "1. What did you like the most on the website?"
$start = $found.Start
$fullText = "4. If applicable, describe any major challenge you may have had during the tasks?"
$whole = $d.Range($found.Start, $found.End)
$whole.Text = $fullText

$pieces = @("4. ", "If applicable, describe any major", " challenge", " you may have had ", "during the tasks?")
$pos = $start
foreach ($piece in $pieces) {
    $len = $piece.Length
    $r = $d.Range($pos, $pos + $len)
    $r.Font.HighlightColorIndex = 4   # wdBrightGreen
    $pos = $pos + $len
}

# ---------------------------------------------------------------------------
# 2) "2. Is there anything you would change about it?"
#    -> "5. How would you rate the colour palette of this website?" (yellow)
# ---------------------------------------------------------------------------
$found = Find-Text $d "2. Is there anything you would change about it?"
$start = $found.Start
$newText = "5. How would you rate the colour palette of this website?"
$whole = $d.Range($found.Start, $found.End)
$whole.Text = $newText
$r = $d.Range($start, $start + $newText.Length)
$r.Font.HighlightColorIndex = 7   # wdYellow

# ---------------------------------------------------------------------------
# 3) "3. On a scale of 10, how likely is it that you would recommend this
#     website to a friend or colleague?"
#    -> "6. Were you able to clearly differentiate between links and other
#        content?" (yellow)
# ---------------------------------------------------------------------------
$found = Find-Text $d "3. On a scale of 10, how likely is it that you would recommend this website to a friend or colleague?"
$start = $found.Start
$newText = "6. Were you able to clearly differentiate between links and other content?"
$whole = $d.Range($found.Start, $found.End)
$whole.Text = $newText
$r = $d.Range($start, $start + $newText.Length)
$r.Font.HighlightColorIndex = 7   # wdYellow

# ---------------------------------------------------------------------------
# 4) "4. Did you face any challenges during the tasks?"
#    -> "7. Was the size of the content, menus, drop-downs, and other
#        features appropriate?" (yellow)
# ---------------------------------------------------------------------------
$found = Find-Text $d "4. Did you face any challenges during the tasks?"
$start = $found.Start
$newText = "7. Was the size of the content, menus, drop-downs, and other features appropriate?"
$whole = $d.Range($found.Start, $found.End)
$whole.Text = $newText
$r = $d.Range($start, $start + $newText.Length)
$r.Font.HighlightColorIndex = 7   # wdYellow

# ---------------------------------------------------------------------------
# 5) Remove the six paragraphs that are no longer part of the questionnaire:
#    "5. How would you rate the colour palette of this website?"
#    "6. Were you able to clearly differentiate between links and other content?"
#    "7. Was the size of the content, menus, drop-downs, and other features appropriate?"
#    "8. How easy was it to navigate through the website?"
#    "9. What are some aspects you disliked?"
#    "10. Could the content have been presented in a better way?"
#
# NOTE: by this point in the script, questions 1-4 were renumbered to 4-7,
# which makes their NEW text identical to some of the OLD question 5/6/7
# paragraphs we are about to delete (e.g. both a kept and a to-be-deleted
# paragraph now read "5. How would you rate the colour palette of this
# website?"). A plain text search would therefore be ambiguous, so the six
# paragraphs to remove are addressed by their (stable, unaffected-by-the-
# above 1-for-1 text substitutions) paragraph index instead.
# ---------------------------------------------------------------------------
$firstDelPara = $d.Paragraphs.Item(15)
$lastDelPara = $d.Paragraphs.Item(20)
Write-Host "Deleting paragraphs:" $firstDelPara.Range.Text "..." $lastDelPara.Range.Text
$delRange = $d.Range($firstDelPara.Range.Start, $lastDelPara.Range.End)
$delRange.Delete()

# ---------------------------------------------------------------------------
# 6) "11. We would love to get your suggestions to help us improve our
#     website " -> split into two runs: the sentence (green highlight, no
#     trailing space) and a trailing space run (no highlight, unchanged).
# ---------------------------------------------------------------------------
$sentence = "11. We would love to get your suggestions to help us improve our website"
$found = Find-Text $d $sentence
$found.Font.HighlightColorIndex = 4   # wdBrightGreen

Write-Host "Edit complete"
